$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: new marks in columns F and G (raises the SUM in N6 from 10 to 20)
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = 5

# Row 14: new mark in column F (raises N14 from 13 to 18)
$ws.Range("F14").Value = 5

# Row 18: new mark in column F (raises N18 from 20 to 25)
$ws.Range("F18").Value = 5

# Row 21: new mark in column F (raises N21 from 5 to 10)
$ws.Range("F21").Value = 5

# Move the active selection to the newly entered G6 cell, matching the
# author's final cursor position / frozen-pane scroll state.
$ws.Range("G6").Select()
